$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: B1:F1
$ws.Range("B1").Value = "Exp 1"
$ws.Range("C1").Value = "Exp 2"
$ws.Range("D1").Value = "Exp 3"
$ws.Range("E1").Value = "Exp 4"
$ws.Range("F1").Value = "Exp 5"

# Copy the existing header style (bold/centered/bordered) from B1 to C1:F1
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2: Accuracy values across B2:F2
$ws.Range("B2").Value = 0.9635036496350365
$ws.Range("C2").Value = 0.9854014598540146
$ws.Range("D2").Value = 0.9854014598540146
$ws.Range("E2").Value = 0.9343065693430657
$ws.Range("F2").Value = 0.9708029197080292

# Row 3: rename label and set Sensitivity values
$ws.Range("A3").Value = "Sensitivity"
$ws.Range("B3").Value = 0.9387755102040817
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.9666666666666667
$ws.Range("E3").Value = 0.8837209302325582
$ws.Range("F3").Value = 0.9574468085106383

# Row 4: rename label and set Geometric Mean values
$ws.Range("A4").Value = "Geometric Mean"
$ws.Range("B4").Value = 0.9578307278188505
$ws.Range("C4").Value = 0.9881652636251156
$ws.Range("D4").Value = 0.983192080250175
$ws.Range("E4").Value = 0.9198455219574726
$ws.Range("F4").Value = 0.9675588936937934

# Remove old rows 5 (Specificity) and 6 (old Geometric Mean)
$ws.Rows("5:6").Delete()
